$d = $word.ActiveDocument
$sec = $d.Sections.First
$hEven = $sec.Headers.Item(2)
$hFirst = $sec.Headers.Item(3)
$fEven = $sec.Footers.Item(2)
$fFirst = $sec.Footers.Item(3)
Write-Output ("hEven.LinkToPrevious before=" + $hEven.LinkToPrevious)
$hEven.LinkToPrevious = $true
$hFirst.LinkToPrevious = $true
$fEven.LinkToPrevious = $true
$fFirst.LinkToPrevious = $true
Write-Output ("hEven.LinkToPrevious after=" + $hEven.LinkToPrevious)
$sec.PageSetup.DifferentFirstPageHeaderFooter = $false
$sec.PageSetup.OddAndEvenPagesHeaderFooter = $false
Write-Output "done"
